$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''43.142.82'
$ws.Range('E2').Value = '  +2.40%  '
$ws.Range('D3').Value = '''2.315.83'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''303.02'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').Value = '''102.22'
$ws.Range('E6').Value = '  +7.21%  '
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +6.16%  '
$ws.Range('D10').Value = '''36.19'
$ws.Range('E10').Value = '  +9.48%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('E12').Value = '  +3.74%  '
$ws.Range('E13').Value = '  +14.69%  '
$ws.Range('E14').Value = '  +3.90%  '
$ws.Range('D15').Value = '''2.671.80'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').Value = '''2.299.34'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').Value = '''0.810'
$ws.Range('E17').Value = '  +3.73%  '
$ws.Range('D18').Value = '''43.052.50'
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('D19').Value = '''12.66'
$ws.Range('E19').Value = '  +8.38%  '
$ws.Range('D20').Value = '''6.20'
$ws.Range('E20').Value = '  +3.98%  '
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').Value = '''67.94'
$ws.Range('D23').Value = '''237.46'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = '''2.22'
$ws.Range('E24').Value = '  +13.15%  '
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '''24.81'
$ws.Range('E27').Value = '  +4.82%  '
$ws.Range('E28').Value = '  +3.79%  '
$ws.Range('D29').Value = '''168.25'
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('D30').Value = '''34.61'
$ws.Range('E30').Value = '  +3.44%  '
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +4.25%  '
$ws.Range('E34').Value = '  +3.59%  '
$ws.Range('D35').Value = '''17.35'
$ws.Range('E35').Value = '  +4.61%  '
$ws.Range('E36').Value = '  +3.52%  '
$ws.Range('D37').Value = '''0.0694'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('E38').Value = '  +4.67%  '
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').Value = '''1.80'
$ws.Range('E40').Value = '  +4.82%  '
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('D43').Value = '''1.990.99'
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('E44').Value = '  +4.70%  '
$ws.Range('D45').Value = '''10.21'
$ws.Range('E45').Value = '  +7.43%  '
$ws.Range('D46').Value = '''17.69'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  +4.83%  '
$ws.Range('D48').Value = '''56.03'
$ws.Range('E48').Value = '  +7.55%  '
$ws.Range('D49').Value = '''2.545.29'
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('E50').Value = '  +4.09%  '
$ws.Range('E51').Value = '  +2.73%  '
